$wb = $excel.ActiveWorkbook

# Map of worksheet name -> hashtable of row number -> new value for column F
$changes = @{
    "展览" = @{ 3 = 3251; 5 = 2324; 8 = 1320; 9 = 1058; 10 = 274; 11 = 490; 16 = 8167; 19 = 231; 20 = 245; 22 = 466; 27 = 1893; 28 = 997; 30 = 1708; 34 = 14; 35 = 7; 36 = 61; 40 = 203; 41 = 368; 43 = 239 }
    "全部类型" = @{ 5 = 3251; 7 = 2324; 10 = 1320; 12 = 1058; 13 = 274; 14 = 490; 17 = 8167; 21 = 231; 22 = 245; 24 = 466; 29 = 1893; 30 = 997; 32 = 1708; 35 = 14; 36 = 7; 37 = 61; 41 = 203; 42 = 368; 49 = 239 }
}

foreach ($sheetName in $changes.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowChanges = $changes[$sheetName]
    foreach ($row in $rowChanges.Keys) {
        $ws.Range("F$row").Value = $rowChanges[$row]
    }
}
